# Baker Andina backup workbook: push every "due date" cell forward by 2 years
# (2024 -> 2026, 2023 -> 2025). The workbook has two sheets, "Inglês" and
# "Português", each laid out as a per-collaborator roadmap whose last populated
# month/year cell in a row is the "due" date that needs to roll forward.

$wb = $excel.ActiveWorkbook

$wsEN = $wb.Worksheets.Item("Inglês")
$wsEN.Range("J2").Value = "05/2026"
$wsEN.Range("M3").Value = "07/2026"
$wsEN.Range("I4").Value = "07/2026"
$wsEN.Range("L5").Value = "04/2026"
$wsEN.Range("F6").Value = "07/2026"
$wsEN.Range("M7").Value = "05/2026"
$wsEN.Range("J8").Value = "05/2026"
$wsEN.Range("G9").Value = "04/2026"
$wsEN.Range("F10").Value = "04/2026"
$wsEN.Range("G11").Value = "07/2026"
$wsEN.Range("K12").Value = "07/2026"
$wsEN.Range("J13").Value = "07/2026"
$wsEN.Range("F14").Value = "07/2026"
$wsEN.Range("I15").Value = "04/2026"
$wsEN.Range("H16").Value = "05/2026"
$wsEN.Range("J17").Value = "07/2026"
$wsEN.Range("H18").Value = "04/2026"
$wsEN.Range("H19").Value = "05/2026"
$wsEN.Range("I20").Value = "04/2026"
$wsEN.Range("H21").Value = "07/2026"
$wsEN.Range("N22").Value = "05/2026"
$wsEN.Range("F23").Value = "07/2026"
$wsEN.Range("H24").Value = "05/2026"
$wsEN.Range("N25").Value = "05/2026"
$wsEN.Range("H26").Value = "07/2026"
$wsEN.Range("N27").Value = "04/2026"
$wsEN.Range("J28").Value = "04/2026"
$wsEN.Range("H29").Value = "04/2026"
$wsEN.Range("I30").Value = "07/2026"
$wsEN.Range("M31").Value = "07/2026"
$wsEN.Range("N32").Value = "05/2026"
$wsEN.Range("F33").Value = "07/2026"
$wsEN.Range("M34").Value = "07/2026"
$wsEN.Range("G35").Value = "04/2026"
$wsEN.Range("H36").Value = "06/2026"
$wsEN.Range("H37").Value = "05/2026"
$wsEN.Range("N38").Value = "06/2026"
$wsEN.Range("N39").Value = "05/2026"
$wsEN.Range("M40").Value = "04/2026"
$wsEN.Range("L41").Value = "04/2026"
$wsEN.Range("J42").Value = "04/2026"
$wsEN.Range("L43").Value = "04/2026"
$wsEN.Range("N44").Value = "05/2026"
$wsEN.Range("O45").Value = "04/2026"
$wsEN.Range("K46").Value = "04/2026"
$wsEN.Range("J47").Value = "04/2026"
$wsEN.Range("K48").Value = "07/2026"
$wsEN.Range("N49").Value = "07/2026"
$wsEN.Range("L50").Value = "07/2026"
$wsEN.Range("M51").Value = "04/2026"
$wsEN.Range("N52").Value = "05/2026"
$wsEN.Range("H53").Value = "04/2026"
$wsEN.Range("I55").Value = "04/2026"
$wsEN.Range("O56").Value = "04/2026"
$wsEN.Range("N57").Value = "05/2026"
$wsEN.Range("H58").Value = "05/2026"
$wsEN.Range("L59").Value = "04/2026"
$wsEN.Range("L60").Value = "04/2026"
$wsEN.Range("F61").Value = "04/2026"
$wsEN.Range("H62").Value = "07/2026"
$wsEN.Range("H63").Value = "07/2026"
$wsEN.Range("J64").Value = "07/2026"
$wsEN.Range("L65").Value = "05/2026"
$wsEN.Range("L66").Value = "04/2026"
$wsEN.Range("H67").Value = "07/2026"
$wsEN.Range("N68").Value = "07/2026"
$wsEN.Range("H69").Value = "04/2026"
$wsEN.Range("F70").Value = "07/2026"
$wsEN.Range("L71").Value = "07/2026"
$wsEN.Range("H72").Value = "07/2026"
$wsEN.Range("O73").Value = "07/2026"
$wsEN.Range("I74").Value = "07/2026"
$wsEN.Range("I75").Value = "05/2026"
$wsEN.Range("H76").Value = "07/2026"
$wsEN.Range("G77").Value = "07/2026"
$wsEN.Range("O78").Value = "12/2025"
$wsEN.Range("L79").Value = "04/2026"
$wsEN.Range("F80").Value = "04/2026"
$wsEN.Range("L81").Value = "07/2026"
$wsEN.Range("K82").Value = "07/2026"
$wsEN.Range("J83").Value = "06/2026"
$wsEN.Range("O84").Value = "04/2026"
$wsEN.Range("I85").Value = "05/2026"
$wsEN.Range("J86").Value = "04/2026"
$wsEN.Range("K87").Value = "07/2026"
$wsEN.Range("M88").Value = "04/2026"
$wsEN.Range("L89").Value = "07/2026"
$wsEN.Range("M90").Value = "07/2026"
$wsEN.Range("M91").Value = "07/2026"
$wsEN.Range("J92").Value = "07/2026"
$wsEN.Range("F93").Value = "04/2026"
$wsEN.Range("H94").Value = "04/2026"
$wsEN.Range("H95").Value = "09/2025"
$wsEN.Range("L96").Value = "07/2026"
$wsEN.Range("H97").Value = "04/2026"
$wsEN.Range("G98").Value = "07/2026"
$wsEN.Range("K99").Value = "04/2026"
$wsEN.Range("L100").Value = "04/2026"
$wsEN.Range("L101").Value = "04/2026"
$wsEN.Range("I102").Value = "07/2026"
$wsEN.Range("L103").Value = "05/2026"
$wsEN.Range("N104").Value = "05/2026"
$wsEN.Range("H105").Value = "07/2026"
$wsEN.Range("N106").Value = "05/2026"
$wsEN.Range("J107").Value = "04/2026"
$wsEN.Range("J108").Value = "04/2026"
$wsEN.Range("F109").Value = "04/2026"
$wsEN.Range("N110").Value = "05/2026"
$wsEN.Range("F112").Value = "04/2026"
$wsEN.Range("K113").Value = "05/2026"
$wsEN.Range("H114").Value = "07/2026"

$wsPT = $wb.Worksheets.Item("Português")
$wsPT.Range("D4").Value = "04/2026"
$wsPT.Range("D5").Value = "07/2026"
$wsPT.Range("D7").Value = "04/2026"
